# Apply crypto price/volume updates scraped on Tue Mar 26 18:35:16 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain TEXT (not auto-converted to a
# number) without leaving a residual style on the cell - format as text,
# assign, then clear the format back off so the cell matches the original
# (unstyled) inline-string cells.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = '70.224.36'
$ws.Range("E2").Value = '  -0.37%  '
# Row 3
$ws.Range("D3").Value = '3.583.23'
$ws.Range("E3").Value = '  -1.16%  '
# Row 4
$ws.Range("E4").Value = '  -0.08%  '
# Row 5
Set-TextValue $ws.Range("D5") '576.71'
$ws.Range("E5").Value = '  -2.75%  '
# Row 6
Set-TextValue $ws.Range("D6") '189.53'
$ws.Range("E6").Value = '  -1.30%  '
# Row 7
$ws.Range("E7").Value = '  -2.26%  '
# Row 8
$ws.Range("D8").Value = '3.578.23'
$ws.Range("E8").Value = '  -0.59%  '
# Row 9
$ws.Range("E9").Value = '  -0.07%  '
# Row 10
$ws.Range("E10").Value = '  -0.75%  '
# Row 11
$ws.Range("E11").Value = '  +0.12%  '
# Row 12
Set-TextValue $ws.Range("D12") '56.00'
$ws.Range("E12").Value = '  -2.89%  '
# Row 13
Set-TextValue $ws.Range("D13") '0.0000303'
$ws.Range("E13").Value = '  +2.38%  '
# Row 14
$ws.Range("E14").Value = '  -0.92%  '
# Row 15
$ws.Range("D15").Value = '4.162.59'
$ws.Range("E15").Value = '  -1.19%  '
# Row 16
Set-TextValue $ws.Range("D16") '19.91'
$ws.Range("E16").Value = '  +3.06%  '
# Row 17
$ws.Range("D17").Value = '3.588.11'
$ws.Range("E17").Value = '  -1.22%  '
# Row 18
$ws.Range("D18").Value = '70.082.09'
# Row 19
Set-TextValue $ws.Range("D19") '12.65'
$ws.Range("E19").Value = '  +0.65%  '
# Row 21
$ws.Range("E21").Value = '  -0.57%  '
# Row 22
Set-TextValue $ws.Range("D22") '477.94'
$ws.Range("E22").Value = '  -3.19%  '
# Row 23
Set-TextValue $ws.Range("D23") '18.99'
$ws.Range("E23").Value = '  +13.81%  '
# Row 24
Set-TextValue $ws.Range("D24") '5.11'
$ws.Range("E24").Value = '  -7.40%  '
# Row 25
$ws.Range("E25").Value = '  -1.80%  '
# Row 26
Set-TextValue $ws.Range("D26") '91.94'
$ws.Range("E26").Value = '  +1.62%  '
# Row 27
$ws.Range("E27").Value = '  -2.10%  '
# Row 28
Set-TextValue $ws.Range("D28") '11.04'
# Row 29
Set-TextValue $ws.Range("D29") '9.35'
$ws.Range("E29").Value = '  -0.07%  '
# Row 30
Set-TextValue $ws.Range("D30") '32.31'
$ws.Range("E30").Value = '  -0.11%  '
# Row 31
Set-TextValue $ws.Range("D31") '7.70'
$ws.Range("E31").Value = '  +1.21%  '
# Row 32
$ws.Range("E32").Value = '  +3.49%  '
# Row 33
Set-TextValue $ws.Range("D33") '12.21'
$ws.Range("E33").Value = '  -0.11%  '
# Row 34
Set-TextValue $ws.Range("D34") '66.56'
$ws.Range("E34").Value = '  +1.89%  '
# Row 35
Set-TextValue $ws.Range("D35") '586.23'
$ws.Range("E35").Value = '  -4.50%  '
# Row 36
Set-TextValue $ws.Range("D36") '39.12'
$ws.Range("E36").Value = '  +3.35%  '
# Row 37
$ws.Range("E37").Value = '  +0.00%  '
# Row 38
$ws.Range("E38").Value = '  -3.88%  '
# Row 39
$ws.Range("E39").Value = '  -1.06%  '
# Row 40
Set-TextValue $ws.Range("D40") '3.29'
$ws.Range("E40").Value = '  +20.28%  '
# Row 41
Set-TextValue $ws.Range("D41") '0.139'
$ws.Range("E41").Value = '  -5.95%  '
# Row 42
Set-TextValue $ws.Range("D42") '3.50'
$ws.Range("E42").Value = '  -5.07%  '
# Row 43
$ws.Range("D43").Value = '3.241.31'
$ws.Range("E43").Value = '  -3.45%  '
# Row 44
Set-TextValue $ws.Range("D44") '2.86'
$ws.Range("E44").Value = '  +7.47%  '
# Row 45
$ws.Range("E45").Value = '  +1.57%  '
# Row 46
$ws.Range("E46").Value = '  -0.26%  '
# Row 47
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D47") '9.46'
$ws.Range("E47").Value = '  +4.27%  '
# Row 48
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D48") '3.34'
$ws.Range("E48").Value = '  -0.46%  '
# Row 49
$ws.Range("E49").Value = '  +0.16%  '
# Row 50
$ws.Range("E50").Value = '  -0.05%  '
# Row 51
Set-TextValue $ws.Range("D51") '3.15'
$ws.Range("E51").Value = '  -4.47%  '
